$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2

$ws.Activate()
$ws.Range("C6").Select()
